$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E3 becomes blank (was -5.7) ---
# ClearContents() alone drops the cell entirely when it is saved back out, so
# nudge the cell's formatting (re-apply its existing style) to make the
# engine keep an (empty) cell record at E3 instead of omitting it.
$ws.Range("E3").ClearContents()
$ws.Range("E3").Style = "Normal"

# --- Row "RM 232" (row 26) is removed entirely; rows below shift up ---
$ws.Rows(26).Delete()

# --- Row "SC 92" (now row 27 after the first delete) is removed too ---
$ws.Rows(27).Delete()

# After the two deletions the remaining "SC *" rows have shifted up by two,
# landing on rows 26-33. A few of those rows also got value edits:

# Row 26 is now "SC 5" -> B26 goes from blank to -20.2
$ws.Range("B26").Value = -20.2

# Row 27 is now "SC 101" -> B27 goes from -20.4 to blank
$ws.Range("B27").ClearContents()
$ws.Range("B27").Style = "Normal"

# Row 33 is now "SC 232" -> B33 and E33 go from blank to real values
$ws.Range("B33").Value = -19.5
$ws.Range("E33").Value = -10.7
